# Nearest greater/smaller to right/left
#
# This adds two new "Stack" problems (and their accompanying metadata) to the
# "Stack" worksheet, makes that worksheet the active tab of the workbook, and
# updates the now-inactive "General problems" worksheet's remembered
# selection.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General problems")
$wsStack   = $wb.Worksheets.Item("Stack")

# --- "General problems" keeps its scroll position but the remembered
#     selection moves on to the next empty row (A17); it also stops being the
#     active/selected tab once "Stack" is activated below. ---
$wsGeneral.Activate()
$wsGeneral.Range("A17").Select()

# --- "Stack" worksheet: append the two new rows ---
$wsStack.Activate()

# Row 6: Nearest greater to right/left
$wsStack.Range("A5").Copy()
$wsStack.Range("A6").PasteSpecial(-4122)   # xlPasteFormats (reuse date style)
$wsStack.Range("A6").Value = 44406
$wsStack.Range("B6").Value = "Nearest greater to right/left"
$wsStack.Range("C6").Value = "Done - D"
$wsStack.Range("D6").Value = "Medium"

$wsGeneral.Range("E15").Copy()
$wsStack.Range("E6").PasteSpecial(-4122)   # xlPasteFormats (reuse wrap-text style)
$wsStack.Range("E6").Value = "NearestGreaterToRight`nNearestGreaterToLeft"

# Row 7: Nearest smaller to right/left
$wsStack.Range("A5").Copy()
$wsStack.Range("A7").PasteSpecial(-4122)
$wsStack.Range("A7").Value = 44407
$wsStack.Range("B7").Value = "Nearest smaller to right/left"
$wsStack.Range("C7").Value = "Done - D"
$wsStack.Range("D7").Value = "Medium"

$wsGeneral.Range("E15").Copy()
$wsStack.Range("E7").PasteSpecial(-4122)
$wsStack.Range("E7").Value = "NearestSmallerToRight`nNearestSmallerToLeft"

# Row 8: a blank placeholder row, only the date-formatted style in column A
$wsStack.Range("A5").Copy()
$wsStack.Range("A8").PasteSpecial(-4122)

# The two wrapped, two-line descriptions make their rows taller.
$wsStack.Rows.Item(6).RowHeight = 29
$wsStack.Rows.Item(7).RowHeight = 29

# "Stack" becomes the active tab/sheet, with F6 remembered as the selection.
$wsStack.Range("F6").Select()
